$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.306.45'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.45%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.495.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.30%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.13%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.44'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.495.13'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.34%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").Value = '  -1.16%  '

$ws.Range("E10").Value = '  -0.05%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.16'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.376'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.31%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.089.87'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.20%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.120'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.68%  '

$ws.Range("E15").Value = '  +0.15%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.496.93'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.01%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.348.53'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -9.73%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.84'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.76%  '

$ws.Range("E20").Value = '  +1.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '388.32'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.93%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.635.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.23%  '

$ws.Range("E24").Value = '  -2.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.99%  '

$ws.Range("E26").Value = '  -0.17%  '

$ws.Range("E27").Value = '  -0.81%  '

$ws.Range("E28").Value = '  +0.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.54'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.25%  '

$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.40'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.39%  '

$ws.Range("E32").Value = '  +0.05%  '

$ws.Range("E33").Value = '  -1.35%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.517.01'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.23%  '

$ws.Range("E35").Value = '  -0.01%  '

$ws.Range("E36").Value = '  +2.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '23.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.73%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.25'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.75%  '

$ws.Range("E39").Value = '  -1.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.53'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.35%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '162.06'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0782'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.805'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.22%  '

$ws.Range("E44").Value = '  +0.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '25.40'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.69%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.82'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.42%  '

$ws.Range("E47").Value = '  +0.13%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.18'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.17%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.66'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.61%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.477.27'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.37%  '

$ws.Range("E51").Value = '  -2.23%  '
